# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Swap the display order of "Santa Lucia" / "Timor Oriental" rows
# - Update Pakistan, Ucrania, Uzbekistan and Tailandia COVID figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Septiembre de 2020 a las 07:15"

# --- Swap Santa Lucia / Timor Oriental labels (rows 206 / 207) ---
$ws.Range("A206").Value = "Santa Lucia"
$ws.Range("A207").Value = "Timor Oriental"

# --- Row 22: Pakistan ---
$ws.Range("B22").Value = 310275
$ws.Range("C22").Value = 694
$ws.Range("D22").Value = 295613
$ws.Range("E22").Value = 8205
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 6457

# --- Row 28: Ucrania ---
$ws.Range("D28").Value = 86873
$ws.Range("E28").Value = 104728

# --- Row 60: Uzbekistan ---
$ws.Range("B60").Value = 54953
$ws.Range("C60").Value = 134
$ws.Range("D60").Value = 51458
$ws.Range("E60").Value = 3041
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 454

# --- Row 141: Tailandia ---
$ws.Range("B141").Value = 3523
$ws.Range("C141").Value = 1
$ws.Range("D141").Value = 3367
$ws.Range("E141").Value = 97
